$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "Save" - reuse the same formatting as the other
# header cells (e.g. G1: bold font, border, centered alignment).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the "Save" column values for each data row (2..14).
$values = @(1, 0, 0, 0, 1, 0, 0, 0, 1, 0, 1, 1, 1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
